$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (supplemental-figures update).
$values = @(
    @(-18.01946042125931,-0.4043409094625655,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,0.7020159779691936,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-0.418785241362327,0.1412825825299654,-18.01946042125931,3.897186273245978,-18.01946042125931,1.211472390416359,-18.01946042125931,1.655443872148824,-18.01946042125931),
    @(-18.01946042125931,0.4305466005793847,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(3.238128265203629,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,1.478486896919308,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(3.401117167206553,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,0.1691718854451046,-18.01946042125931,2.168975093251183),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,2.103874176637821,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,1.325574418495182),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,2.079287494305293,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,1.293238658792612,1.90574773256668),
    @(-18.01946042125931,-18.01946042125931,0.9854649996562829,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,2.42337477285154),
    @(-18.01946042125931,-18.01946042125931,-0.342449695318466,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,2.227814641633442,-18.01946042125931),
    @(-18.01946042125931,0.1517799151992887,-0.2512179328067725,-18.01946042125931,-18.01946042125931,-18.01946042125931,1.813913042625469,0.07251723746267358,2.019234897393729,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,-18.01946042125931,2.53605084069446,0.2617497081430605,2.501067144249472,-18.01946042125931),
    @(-18.01946042125931,-18.01946042125931,3.008451648064503,-18.01946042125931,-18.01946042125931,-18.01946042125931,1.625383408794287,0.8413264980532424,-18.01946042125931,-18.01946042125931),
    @(-18.01946042125931,3.254765858244162,2.856232693758332,-18.01946042125931,2.350636283485992,-18.01946042125931,1.426065735876561,3.72318495250233,-18.01946042125931,1.95802774455155),
    @(-18.01946042125931,2.698772252306265,-18.01946042125931,3.119151995586959,-18.01946042125931,4.321922936655875,1.376234780051335,-18.01946042125931,-18.01946042125931,-18.01946042125931)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
